$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Consumo" (column C) measurement values for each row.
$ws.Range("C3").Value = 1575
$ws.Range("C4").Value = 1400
$ws.Range("C5").Value = 855
$ws.Range("C6").Value = 286
$ws.Range("C7").Value = 148
$ws.Range("C8").Value = 4562
$ws.Range("C9").Value = 454

# Column F ("periodoDeImputacion") moves from a plain "2022" placeholder to
# real dates, formatted as month-year.
$ws.Range("F3:F9").NumberFormat = "mmm-yy"
$ws.Range("F3").Value = 44562
$ws.Range("F4").Value = 44593
$ws.Range("F5").Value = 44682
$ws.Range("F6").Value = 44652
$ws.Range("F7").Value = 44652
$ws.Range("F8").Value = 44197
$ws.Range("F9").Value = 44531

# Rows with wrapped "Tipo de Consumo" text end up taller once the sheet is
# recalculated/redisplayed by Excel.
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
